$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper fragments (flat-OPC wrapped WordprocessingML) used with
# Range.InsertXML to get exact run/bookmark/pPr structure that plain
# Range.Text / InsertBefore text APIs cannot produce (those always merge
# into a single run).
# ---------------------------------------------------------------------------
$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgClose = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ===========================================================================
# 4) Append a new, empty "List Paragraph" bullet (numId 1) right after the
#    last paragraph ("Select input using same name style as radio button.").
#    Doing this first (it is the very end of the document) means none of the
#    earlier paragraph indices used below are disturbed by it.
# ===========================================================================
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$newBulletXml = $pkgOpen + '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p></w:body>' + $pkgClose
$endPoint.InsertXML($newBulletXml)

# ===========================================================================
# 3) "[select][Name]" -> "[select]" + " " + "_[" + "TabName" + "]_" + "[Name]"
#    (split into separate runs, TabName wrapped in spell-check markers, no
#    bookmark here).
# ===========================================================================
$pSelect = $d.Paragraphs.Item(21)
if ($pSelect.Range.Text -notmatch "\[select\]\[Name\]") {
    throw "Paragraph 21 did not contain expected '[select][Name]' text: " + $pSelect.Range.Text
}
$rSelect = $d.Range($pSelect.Range.Start, $pSelect.Range.End - 1)
$selectXml = $pkgOpen + '<w:body><w:p>' + `
    '<w:r><w:t>[select]</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>_[</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>TabName</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>]_</w:t></w:r>' + `
    '<w:r><w:t>[Name]</w:t></w:r>' + `
    '</w:p></w:body>' + $pkgClose
$rSelect.InsertXML($selectXml)

# ===========================================================================
# 2) "[is][Name]" -> "[is]" + " " + "_[" + "TabName" + "]_" + bookmark + "[Name]"
#    (same split as above, but with the _GoBack bookmark reinstated right
#    before the final "[Name]" run).
# ===========================================================================
$pIs = $d.Paragraphs.Item(17)
if ($pIs.Range.Text -notmatch "\[is\]\[Name\]") {
    throw "Paragraph 17 did not contain expected '[is][Name]' text: " + $pIs.Range.Text
}
$rIs = $d.Range($pIs.Range.Start, $pIs.Range.End - 1)
$isXml = $pkgOpen + '<w:body><w:p>' + `
    '<w:r><w:t>[is]</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>_[</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>TabName</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>]_</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>[Name]</w:t></w:r>' + `
    '</w:p></w:body>' + $pkgClose
$rIs.InsertXML($isXml)

# ===========================================================================
# 1) Remove the stray <w:bookmarkStart/><w:bookmarkEnd/> ("_GoBack") that
#    trails "download_[TabName]_[ButtonName]". Because this engine keeps
#    _GoBack "sticky" to its paragraph (it re-anchors around any in-place
#    text replacement within that paragraph), the only reliable way to drop
#    it is to remove the whole paragraph -- bookmark included -- and
#    recreate a fresh paragraph with the same text in its place.
# ===========================================================================
$pDownload = $d.Paragraphs.Item(13)
if ($pDownload.Range.Text -notmatch "download_\[TabName\]_\[ButtonName\]") {
    throw "Paragraph 13 did not contain expected download button text: " + $pDownload.Range.Text
}
$pDownload.Range.InsertParagraphBefore() | Out-Null
$freshEmptyPara = $d.Paragraphs.Item(13)
$freshEmptyPara.Range.InsertBefore("download_[TabName]_[ButtonName]")
$oldParaWithBookmark = $d.Paragraphs.Item(14)
$d.Range($oldParaWithBookmark.Range.Start, $oldParaWithBookmark.Range.End).Delete()
